$wb = $excel.ActiveWorkbook

# --- Update the AntennaMetadata "SiteName" column for the stationary antennas:
# drop the " Stationary Antenna" suffix so the values read as plain site names.
$ws = $wb.Worksheets.Item("AntennaMetadata")
$ws.Range("B7").Value  = "Red Barn"
$ws.Range("B8").Value  = "Red Barn"
$ws.Range("B9").Value  = "Hitching Post"
$ws.Range("B10").Value = "Hitching Post"
$ws.Range("B11").Value = "Confluence"
$ws.Range("B12").Value = "Confluence"
$ws.Range("B13").Value = "Connectivity Channel Downstream"
$ws.Range("B14").Value = "Connectivity Channel Downstream"
$ws.Range("B15").Value = "Connectivity Channel Side Channel"
$ws.Range("B16").Value = "Connectivity Channel Side Channel"
$ws.Range("B17").Value = "Connectivity Channel Upstream"
$ws.Range("B18").Value = "Connectivity Channel Upstream"

# --- TestTags keeps its own last selection (B10) but is no longer the active tab.
$wsTags = $wb.Worksheets.Item("TestTags")
$wsTags.Activate()
$wsTags.Range("B10").Select()

# --- Re-ran script: re-select AntennaMetadata as the active sheet/tab, with
# B17 as the active cell (matches the refreshed selection state).
$ws.Activate()
$ws.Range("B17").Select()
